$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet text with new conversion rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.57 = 50240.2 pesos`n✅ 50240.2 pesos = 12.5 = 965.44 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the "tasas" sheet numeric rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 79.56
$wsTasas.Range("O10").Value = 3997.11
$wsTasas.Range("N12").Value = 4020
$wsTasas.Range("O12").Value = 77.25
